{"js": "// Update the division-fact worksheet numbers (table cell contents only).\n// The mapping below lists every \"old text -> new text\" change in the\n// exact order the originals appear in the document. Some \"old\" values\n// repeat (e.g. \"44\u00f79=\"), so we resolve each `search()` hit by document\n// order and consume the matching replacement for that occurrence.\nconst replacements = [\n  [\"34\u00f74=\", \"39\u00f74=\"],\n  [\"26\u00f77=\", \"92\u00f77=\"],\n  [\"15\u00f76=\", \"45\u00f76=\"],\n  [\"15\u00f74=\", \"42\u00f73=\"],\n  [\"53\u00f77=\", \"86\u00f73=\"],\n  [\"30\u00f75=\", \"26\u00f74=\"],\n  [\"11\u00f74=\", \"19\u00f78=\"],\n  [\"14\u00f75=\", \"21\u00f79=\"],\n  [\"96\u00f75=\", \"39\u00f75=\"],\n  [\"25\u00f79=\", \"40\u00f76=\"],\n  [\"49\u00f74=\", \"81\u00f76=\"],\n  [\"71\u00f78=\", \"73\u00f77=\"],\n  [\"18\u00f77=\", \"30\u00f79=\"],\n  [\"18\u00f73=\", \"93\u00f76=\"],\n  [\"36\u00f72=\", \"75\u00f79=\"],\n  [\"72\u00f76=\", \"49\u00f76=\"],\n  [\"11\u00f75=\", \"83\u00f72=\"],\n  [\"88\u00f72=\", \"95\u00f74=\"],\n  [\"61\u00f72=\", \"78\u00f72=\"],\n  [\"44\u00f79=\", \"74\u00f73=\"],\n  [\"75\u00f76=\", \"56\u00f79=\"],\n  [\"46\u00f79=\", \"85\u00f77=\"],\n  [\"44\u00f79=\", \"57\u00f76=\"],\n  [\"30\u00f74=\", \"16\u00f77=\"],\n  [\"45\u00f74=\", \"99\u00f75=\"],\n];\n\n// Group the target replacement texts by their source text, preserving\n// the order in which each duplicate source should be consumed.\nconst queues = new Map();\nfor (const [from, to] of replacements) {\n  if (!queues.has(from)) queues.set(from, []);\n  queues.get(from).push(to);\n}\n\nconst body = context.document.body;\nconst searchResultsBySource = new Map();\nfor (const from of queues.keys()) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  searchResultsBySource.set(from, results);\n}\nawait context.sync();\n\nfor (const [from, queue] of queues) {\n  const results = searchResultsBySource.get(from);\n  if (results.items.length !== queue.length) {\n    throw new Error(\n      `Expected ${queue.length} occurrence(s) of \"${from}\" but found ${results.items.length}`\n    );\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(queue[i], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the division-fact worksheet numbers (table cell contents only).\n# The list below gives every \"old text -> new text\" change in the exact\n# order the originals appear in the document. A few \"old\" values repeat\n# (e.g. \"44\u00f79=\"); re-using the same Range/Find object for every call makes\n# each Execute() resume searching just after the previous match, so\n# repeated source text is resolved strictly in document order.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"34\u00f74=\", \"39\u00f74=\"),\n    @(\"26\u00f77=\", \"92\u00f77=\"),\n    @(\"15\u00f76=\", \"45\u00f76=\"),\n    @(\"15\u00f74=\", \"42\u00f73=\"),\n    @(\"53\u00f77=\", \"86\u00f73=\"),\n    @(\"30\u00f75=\", \"26\u00f74=\"),\n    @(\"11\u00f74=\", \"19\u00f78=\"),\n    @(\"14\u00f75=\", \"21\u00f79=\"),\n    @(\"96\u00f75=\", \"39\u00f75=\"),\n    @(\"25\u00f79=\", \"40\u00f76=\"),\n    @(\"49\u00f74=\", \"81\u00f76=\"),\n    @(\"71\u00f78=\", \"73\u00f77=\"),\n    @(\"18\u00f77=\", \"30\u00f79=\"),\n    @(\"18\u00f73=\", \"93\u00f76=\"),\n    @(\"36\u00f72=\", \"75\u00f79=\"),\n    @(\"72\u00f76=\", \"49\u00f76=\"),\n    @(\"11\u00f75=\", \"83\u00f72=\"),\n    @(\"88\u00f72=\", \"95\u00f74=\"),\n    @(\"61\u00f72=\", \"78\u00f72=\"),\n    @(\"44\u00f79=\", \"74\u00f73=\"),\n    @(\"75\u00f76=\", \"56\u00f79=\"),\n    @(\"46\u00f79=\", \"85\u00f77=\"),\n    @(\"44\u00f79=\", \"57\u00f76=\"),\n    @(\"30\u00f74=\", \"16\u00f77=\"),\n    @(\"45\u00f74=\", \"99\u00f75=\")\n)\n\n$range = $d.Range()\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Could not find occurrence of '$oldText'\"\n    }\n}\n"}
